# Update "想去人数" (F column) figures and "最低票价" (G column) sold-out
# status text across the "展览" (sheet 1) and "全部类型" (sheet 4) sheets,
# matching the freshly re-generated output data.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# ---- 展览 (Exhibition) sheet ----
$wsExhibit.Range("F2").Value = 0
$wsExhibit.Range("F3").Value = 0
$wsExhibit.Range("F4").Value = 19150
$wsExhibit.Range("F5").Value = 747
$wsExhibit.Range("G5").Value = "已售罄"
$wsExhibit.Range("F8").Value = 1
$wsExhibit.Range("F9").Value = 7255
$wsExhibit.Range("F10").Value = 462
$wsExhibit.Range("F11").Value = 709
$wsExhibit.Range("F13").Value = 24
$wsExhibit.Range("F14").Value = 0
$wsExhibit.Range("F15").Value = 0
$wsExhibit.Range("F17").Value = 0
$wsExhibit.Range("F18").Value = 0
$wsExhibit.Range("F20").Value = 64
$wsExhibit.Range("F22").Value = 42
$wsExhibit.Range("F24").Value = 47
$wsExhibit.Range("F25").Value = 0
$wsExhibit.Range("F26").Value = 1052
$wsExhibit.Range("F27").Value = 13
$wsExhibit.Range("F28").Value = 147
$wsExhibit.Range("F29").Value = 5216
$wsExhibit.Range("F31").Value = 0
$wsExhibit.Range("F32").Value = 0
$wsExhibit.Range("F33").Value = 22
$wsExhibit.Range("F35").Value = 12357
$wsExhibit.Range("F36").Value = 1305
$wsExhibit.Range("F37").Value = 37
$wsExhibit.Range("F39").Value = 243

# ---- 全部类型 (All types) sheet ----
$wsAll.Range("F3").Value = 1365
$wsAll.Range("G5").Value = "已售罄"
$wsAll.Range("F6").Value = 290
$wsAll.Range("F7").Value = 1083
$wsAll.Range("F9").Value = 7255
$wsAll.Range("F10").Value = 462
$wsAll.Range("F12").Value = 224
$wsAll.Range("F13").Value = 0
$wsAll.Range("F14").Value = 0
$wsAll.Range("F15").Value = 0
$wsAll.Range("F16").Value = 0
$wsAll.Range("F17").Value = 0
$wsAll.Range("F18").Value = 0
$wsAll.Range("F19").Value = 330
$wsAll.Range("F20").Value = 64
$wsAll.Range("F22").Value = 42
$wsAll.Range("F25").Value = 297
$wsAll.Range("F32").Value = 32
$wsAll.Range("F33").Value = 32
$wsAll.Range("F34").Value = 0
$wsAll.Range("F35").Value = 22
$wsAll.Range("F36").Value = 82
$wsAll.Range("F37").Value = 12357
$wsAll.Range("F41").Value = 243
$wsAll.Range("F42").Value = 0
$wsAll.Range("F44").Value = 314
$wsAll.Range("F45").Value = 94
